$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = '26.068.48'
$ws.Range("E2").Value = '  -0.52%  '

$ws.Range("D3").Value = '1.651.20'
$ws.Range("E3").Value = '  -0.47%  '

$ws.Range("E4").Value = '  -0.24%  '

$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = '217.13'
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = '  -0.03%  '

$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = '0.5257'
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = '  +1.87%  '

$ws.Range("E7").Value = '  -0.18%  '

$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = '0.06320'
$ws.Range("D9").Style = "Normal"
$ws.Range("E9").Value = '  +0.79%  '

$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = '20.33'
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = '  -2.07%  '

$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = '0.07794'
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = '  +0.57%  '

$ws.Range("B12").Value = 'WrappedEther'
$ws.Range("C12").Value = 'https://coinranking.com/coin/Mtfb0obXVh59u+wrappedether-weth'
$ws.Range("D12").Value = '1.677.59'
$ws.Range("E12").Value = '  +1.48%  '

$ws.Range("B13").Value = 'Polkadot'
$ws.Range("C13").Value = 'https://coinranking.com/coin/25W7FG7om+polkadot-dot'
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = '4.499'
$ws.Range("D13").Style = "Normal"
$ws.Range("E13").Value = '  +0.32%  '

$ws.Range("D14").Value = '1.878.42'
$ws.Range("E14").Value = '  -0.41%  '

$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = '0.5486'
$ws.Range("D15").Style = "Normal"
$ws.Range("E15").Value = '  +0.59%  '

$ws.Range("D16").Value = '0.0₅8198'
$ws.Range("E16").Value = '  +0.85%  '

$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = '65.43'
$ws.Range("D17").Style = "Normal"
$ws.Range("E17").Value = '  +0.88%  '

$ws.Range("D18").Value = '26.069.13'
$ws.Range("E18").Value = '  -0.56%  '

$ws.Range("E19").Value = '  -0.22%  '

$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = '4.570'
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = '  -0.85%  '

$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = '190.57'
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = '  -0.75%  '

$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = '10.06'
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = '  -0.22%  '

$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = '6.016'
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = '  +0.52%  '

$ws.Range("E24").Value = '  -0.27%  '

$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = '143.64'
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = '  +2.97%  '

$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = '0.1235'
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Value = '  +1.25%  '

$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = '7.222'
$ws.Range("D27").Style = "Normal"

$ws.Range("E28").Value = '  -0.41%  '

$ws.Range("E29").Value = '  -0.87%  '

$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = '0.05831'
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").Value = '  -1.50%  '

$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = '1.271'
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Value = '  -0.14%  '

$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = '3.547'
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = '  +0.04%  '

$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = '3.263'
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value = '  -0.35%  '

$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = '1.581'
$ws.Range("D34").Style = "Normal"
$ws.Range("E34").Value = '  -0.06%  '

$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = '0.9463'
$ws.Range("D35").Style = "Normal"
$ws.Range("E35").Value = '  -1.45%  '

$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = '2.412'
$ws.Range("D36").Style = "Normal"
$ws.Range("E36").Value = '  -0.46%  '

$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = '2.779'
$ws.Range("D37").Style = "Normal"
$ws.Range("E37").Value = '  +0.32%  '

$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = '0.5739'
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = '  +1.37%  '

$ws.Range("E39").Value = '  +1.30%  '

$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = '5.740'
$ws.Range("D40").Style = "Normal"

$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = '0.8421'
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = '  -1.55%  '

$ws.Range("E42").Value = '  -0.11%  '

$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = '103.91'
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = '  +3.24%  '

$ws.Range("D44").Value = '1.029.97'
$ws.Range("E44").Value = '  +1.78%  '

$ws.Range("D45").Value = '1.794.80'
$ws.Range("E45").Value = '  -0.33%  '

$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = '56.95'
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = '  +0.89%  '

$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = '1.003'
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = '  +0.12%  '

$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = '0.4322'
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = '  +2.56%  '

$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = '7.872'
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = '  -2.12%  '

$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = '0.05144'
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = '  -0.44%  '

$ws.Range("E51").Value = '  +1.02%  '
